$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.7321483333333333
$ws.Range("H2").Value = 2.196445
$ws.Range("I2").Value = 0.05113520435363902
$ws.Range("J2").Value = 0.05113520435363902
$ws.Range("M2").Value = 6.045145666666667
$ws.Range("N2").Value = 18.135437
$ws.Range("O2").Value = 0.8160840232643366
$ws.Range("P2").Value = 0.8160840232643367
$ws.Range("Q2").Value = 4.425943324607222
$ws.Range("R2").Value = 39.833489921465
$ws.Range("S2").Value = 0.04173062329936175
$ws.Range("T2").Value = 0.04173062329936176
$ws.Range("G3").Value = 0.7321483333333333
$ws.Range("H3").Value = 2.196445
$ws.Range("I3").Value = 0.05113520435363902
$ws.Range("J3").Value = 0.05113520435363902
$ws.Range("O3").Value = 0.09212864864242169
$ws.Range("P3").Value = 0.09212864864242169
$ws.Range("Q3").Value = 0.4996497490944445
$ws.Range("R3").Value = 4.49684774185
$ws.Range("S3").Value = 0.004711017275154841
$ws.Range("T3").Value = 0.004711017275154841
$ws.Range("G4").Value = 0.7321483333333333
$ws.Range("H4").Value = 2.196445
$ws.Range("I4").Value = 0.05113520435363902
$ws.Range("J4").Value = 0.05113520435363902
$ws.Range("M4").Value = 0.6799149999999999
$ws.Range("N4").Value = 2.039745
$ws.Range("O4").Value = 0.09178732809324164
$ws.Range("P4").Value = 0.09178732809324165
$ws.Range("Q4").Value = 0.4977986340583332
$ws.Range("R4").Value = 4.480187706524999
$ws.Range("S4").Value = 0.004693563779122422
$ws.Range("T4").Value = 0.004693563779122423
$ws.Range("I5").Value = 0.7165747117895102
$ws.Range("J5").Value = 0.7165747117895102
$ws.Range("M5").Value = 6.045145666666667
$ws.Range("N5").Value = 18.135437
$ws.Range("O5").Value = 0.8160840232643366
$ws.Range("P5").Value = 0.8160840232643367
$ws.Range("Q5").Value = 62.02222328659623
$ws.Range("R5").Value = 558.2000095793661
$ws.Range("S5").Value = 0.5847851737666659
$ws.Range("T5").Value = 0.5847851737666659
$ws.Range("I6").Value = 0.7165747117895102
$ws.Range("J6").Value = 0.7165747117895102
$ws.Range("O6").Value = 0.09212864864242169
$ws.Range("P6").Value = 0.09212864864242169
$ws.Range("S6").Value = 0.06601705984850037
$ws.Range("T6").Value = 0.06601705984850037
$ws.Range("I7").Value = 0.7165747117895102
$ws.Range("J7").Value = 0.7165747117895102
$ws.Range("M7").Value = 0.6799149999999999
$ws.Range("N7").Value = 2.039745
$ws.Range("O7").Value = 0.09178732809324164
$ws.Range("P7").Value = 0.09178732809324165
$ws.Range("Q7").Value = 6.975818660323333
$ws.Range("R7").Value = 62.78236794291
$ws.Range("S7").Value = 0.06577247817434384
$ws.Range("T7").Value = 0.06577247817434384
$ws.Range("G8").Value = 2.568000333333333
$ws.Range("H8").Value = 7.704001
$ws.Range("I8").Value = 0.1793560346266988
$ws.Range("J8").Value = 0.1793560346266988
$ws.Range("M8").Value = 6.045145666666667
$ws.Range("N8").Value = 18.135437
$ws.Range("O8").Value = 0.8160840232643366
$ws.Range("P8").Value = 0.8160840232643367
$ws.Range("Q8").Value = 15.52393608704856
$ws.Range("R8").Value = 139.715424783437
$ws.Range("S8").Value = 0.146369594334894
$ws.Range("T8").Value = 0.146369594334894
$ws.Range("G9").Value = 2.568000333333333
$ws.Range("H9").Value = 7.704001
$ws.Range("I9").Value = 0.1793560346266988
$ws.Range("J9").Value = 0.1793560346266988
$ws.Range("O9").Value = 0.09212864864242169
$ws.Range("P9").Value = 0.09212864864242169
$ws.Range("Q9").Value = 1.752514707481111
$ws.Range("R9").Value = 15.77263236733
$ws.Range("S9").Value = 0.01652382909602115
$ws.Range("T9").Value = 0.01652382909602115
$ws.Range("G10").Value = 2.568000333333333
$ws.Range("H10").Value = 7.704001
$ws.Range("I10").Value = 0.1793560346266988
$ws.Range("J10").Value = 0.1793560346266988
$ws.Range("M10").Value = 0.6799149999999999
$ws.Range("N10").Value = 2.039745
$ws.Range("O10").Value = 0.09178732809324164
$ws.Range("P10").Value = 0.09178732809324165
$ws.Range("Q10").Value = 1.746021946638333
$ws.Range("R10").Value = 15.714197519745
$ws.Range("S10").Value = 0.01646261119578361
$ws.Range("T10").Value = 0.01646261119578361
$ws.Range("G11").Value = 0.7579039999999999
$ws.Range("H11").Value = 2.273712
$ws.Range("I11").Value = 0.05293404923015203
$ws.Range("J11").Value = 0.05293404923015203
$ws.Range("M11").Value = 6.045145666666667
$ws.Range("N11").Value = 18.135437
$ws.Range("O11").Value = 0.8160840232643366
$ws.Range("P11").Value = 0.8160840232643367
$ws.Range("Q11").Value = 4.581640081349334
$ws.Range("R11").Value = 41.234760732144
$ws.Range("S11").Value = 0.04319863186341492
$ws.Range("T11").Value = 0.04319863186341493
$ws.Range("G12").Value = 0.7579039999999999
$ws.Range("H12").Value = 2.273712
$ws.Range("I12").Value = 0.05293404923015203
$ws.Range("J12").Value = 0.05293404923015203
$ws.Range("O12").Value = 0.09212864864242169
$ws.Range("P12").Value = 0.09212864864242169
$ws.Range("Q12").Value = 0.5172265321066667
$ws.Range("R12").Value = 4.65503878896
$ws.Range("S12").Value = 0.004876742422745328
$ws.Range("T12").Value = 0.004876742422745328
$ws.Range("G13").Value = 0.7579039999999999
$ws.Range("H13").Value = 2.273712
$ws.Range("I13").Value = 0.05293404923015203
$ws.Range("J13").Value = 0.05293404923015203
$ws.Range("M13").Value = 0.6799149999999999
$ws.Range("N13").Value = 2.039745
$ws.Range("O13").Value = 0.09178732809324164
$ws.Range("P13").Value = 0.09178732809324165
$ws.Range("Q13").Value = 0.5153102981599998
$ws.Range("R13").Value = 4.637792683439999
$ws.Range("S13").Value = 0.004858674943991769
$ws.Range("T13").Value = 0.00485867494399177
